# Update the "Förändrad" date column (C) for rows 2-27 from 2023-10-25 (45224)
# to 2023-11-03 (45233), matching the underlying Excel serial date number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
